$wb = $excel.ActiveWorkbook

# Values are stored as plain text (inline strings) in the source sheet, e.g.
# "4.00", not the number 4. A leading apostrophe forces Excel to keep the
# entry as text instead of auto-converting the numeric-looking string.

# Sheet "部门情况202401" — column N (当年新增普惠户数)
$ws1 = $wb.Worksheets.Item("部门情况202401")
$ws1.Range("N2").Value = "'4.00"
$ws1.Range("N3").Value = "'0.00"
$ws1.Range("N4").Value = "'15.00"
$ws1.Range("N5").Value = "'0.00"
$ws1.Range("N7").Value = "'0.00"

# Sheet "对公业务台账202401" — columns S (三个月小微发生数) and T (三个月发生数)
$ws3 = $wb.Worksheets.Item("对公业务台账202401")
$ws3.Range("S2").Value = "'29.00"
$ws3.Range("T2").Value = "'42.00"
$ws3.Range("S3").Value = "'10.00"
$ws3.Range("T3").Value = "'11.00"
$ws3.Range("S4").Value = "'1.00"
$ws3.Range("T4").Value = "'3.00"
$ws3.Range("S5").Value = "'24.00"
$ws3.Range("T5").Value = "'24.00"
$ws3.Range("S6").Value = "'0.00"
$ws3.Range("T6").Value = "'0.00"
$ws3.Range("S7").Value = "'117.00"
$ws3.Range("T7").Value = "'148.00"
$ws3.Range("S8").Value = "'0.00"
$ws3.Range("T8").Value = "'0.00"
$ws3.Range("S9").Value = "'66.00"
$ws3.Range("T9").Value = "'66.00"
$ws3.Range("S12").Value = "'250.00"
$ws3.Range("T12").Value = "'297.00"
